$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Id (A), Ost (Q) and Nord (R) values between row 5 and row 6
$ws.Range("A5").Value = 111934066
$ws.Range("A6").Value = 111934059

$ws.Range("Q5").Value = 413590.3038565172
$ws.Range("Q6").Value = 413639.6308819132

$ws.Range("R5").Value = 6586912.201658082
$ws.Range("R6").Value = 6586793.951973591

# Move the "Publik kommentar" text from row 5 to row 6
$ws.Range("AC5").ClearContents()
$ws.Range("AC6").Value = "Rätt riklig längs stigen"
